$d = $word.ActiveDocument

function Pin-Boundary($pos) {
    # Touch-and-revert a formatting property on the single character that
    # starts at document position $pos, forcing Word to keep it (and thus
    # the run split at $pos) as its own run instead of silently
    # re-coalescing it with a textually-identical-format neighbour the
    # next time any edit happens in this paragraph.
    $r = $d.Range($pos, $pos + 1)
    $r.Bold = 1
    $r.Bold = 0
}

function Replace-CharAndSplit($searchText, $offset, $newChar, $extraBoundaryOffsets) {
    # Locate the target text anywhere in the document body.
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if (-not $rng.Find.Found) {
        throw "Could not find text: $searchText"
    }
    $segStart = $rng.Start

    # Replace the single changed character. (This can coalesce the whole
    # paragraph's runs, including pre-existing splits - those are
    # re-pinned below.)
    $charStart = $segStart + $offset
    $one = $d.Range($charStart, $charStart + 1)
    $one.Text = $newChar

    # Re-establish the changed character as its own run, separate from
    # both its prefix and its suffix neighbours.
    Pin-Boundary ($segStart + $offset)

    # Re-establish any other pre-existing boundaries in the same
    # paragraph that must stay separate (given as offsets from
    # $segStart, may be negative to reach before the found text).
    foreach ($extra in $extraBoundaryOffsets) {
        Pin-Boundary ($segStart + $extra)
    }
}

# add R3, R1, R2 row is unchanged.

# lw R2, 4(R1)  ->  lw R2, 4(R3)
# Paragraph runs before edit: "lw" | " " | "R2, 4(R1)" - preserve the
# boundary right before "R2, 4(R1)" (offset -1 relative to its start).
Replace-CharAndSplit "R2, 4(R1)" 7 "3" @(-1)

# 100011_00001_00010_00000_00000_000100 -> 100011_00011_00010_00000_00000_000100
# Single-run paragraph, no other boundaries to preserve.
Replace-CharAndSplit "100011_00001_00010_00000_00000_000100" 10 "1" @()

# 8C220004 -> 8C620004 (plain single-run replace, no split needed)
$d.Content.Find.Execute("8C220004", $true, $false, $false, $false, $false, $true, 1, $false, "8C620004", 2) | Out-Null

# 2 <= M[4+R1] -> 2 <= M[4+R3]
# Paragraph runs before edit: "R" | "2 <= M[4+R1]" - preserve the
# boundary right before "2 <= M[4+R1]" (offset -1 relative to its start).
Replace-CharAndSplit "2 <= M[4+R1]" 10 "3" @(-1)

# 10001100001000100000000000000100 -> 10001100011000100000000000000100
# Single-run paragraph, no other boundaries to preserve.
Replace-CharAndSplit "10001100001000100000000000000100" 9 "1" @()
